$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.6167434528639912
$ws.Range("J2").Value = 0.6167434528639911
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1148133333333333
$ws.Range("N2").Value = 0.34444
$ws.Range("O2").Value = 0.03343792635928704
$ws.Range("P2").Value = 0.03343792635928704
$ws.Range("Q2").Value = 0.04551962128444444
$ws.Range("R2").Value = 0.40967659156
$ws.Range("S2").Value = 0.02062262215943856
$ws.Range("T2").Value = 0.02062262215943856

# Row 3
$ws.Range("I3").Value = 0.6167434528639912
$ws.Range("J3").Value = 0.6167434528639911
$ws.Range("O3").Value = 0.9249645515654102
$ws.Range("P3").Value = 0.9249645515654102
$ws.Range("S3").Value = 0.5704658313092443
$ws.Range("T3").Value = 0.5704658313092443

# Row 4
$ws.Range("I4").Value = 0.6167434528639912
$ws.Range("J4").Value = 0.6167434528639911
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1428303333333333
$ws.Range("N4").Value = 0.428491
$ws.Range("O4").Value = 0.04159752207530271
$ws.Range("P4").Value = 0.04159752207530271
$ws.Range("Q4").Value = 0.05662741854544445
$ws.Range("R4").Value = 0.5096467669089999
$ws.Range("S4").Value = 0.02565499939530829
$ws.Range("T4").Value = 0.02565499939530828

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.246372
$ws.Range("H5").Value = 0.739116
$ws.Range("I5").Value = 0.3832565471360088
$ws.Range("J5").Value = 0.3832565471360088
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1148133333333333
$ws.Range("N5").Value = 0.34444
$ws.Range("O5").Value = 0.03343792635928704
$ws.Range("P5").Value = 0.03343792635928704
$ws.Range("Q5").Value = 0.02828679056
$ws.Range("R5").Value = 0.25458111504
$ws.Range("S5").Value = 0.01281530419984849
$ws.Range("T5").Value = 0.01281530419984849

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.246372
$ws.Range("H6").Value = 0.739116
$ws.Range("I6").Value = 0.3832565471360088
$ws.Range("J6").Value = 0.3832565471360088
$ws.Range("O6").Value = 0.9249645515654102
$ws.Range("P6").Value = 0.9249645515654102
$ws.Range("Q6").Value = 0.7824731194280001
$ws.Range("R6").Value = 7.042258074852001
$ws.Range("S6").Value = 0.3544987202561659
$ws.Range("T6").Value = 0.3544987202561659

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.246372
$ws.Range("H7").Value = 0.739116
$ws.Range("I7").Value = 0.3832565471360088
$ws.Range("J7").Value = 0.3832565471360088
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1428303333333333
$ws.Range("N7").Value = 0.428491
$ws.Range("O7").Value = 0.04159752207530271
$ws.Range("P7").Value = 0.04159752207530271
$ws.Range("Q7").Value = 0.035189394884
$ws.Range("R7").Value = 0.316704553956
$ws.Range("S7").Value = 0.01594252267999442
$ws.Range("T7").Value = 0.01594252267999442
